$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:F1) - text swapped due to underlying data reordering
$header = @("L", "Water", "H2", "CO2", "HOAc", "CO")
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(1, $c).Value = $header[$c - 1]
}

# Data rows (A2:F21) - full replacement with re-ordered / expanded dataset
$data = @(
    @(0, -1, 116738344.7276593, 7561327.622822704, -1, 20834017.68103584),
    @(0.12, -1, 108719233.339913, 16282164.97813578, -1, 18987260.75799712),
    @(1.01, -1, 116738344.7276593, 17787591.33748112, -1, 20834017.68103584),
    @(1.04, -1, 116738344.7276593, 16209052.55288416, -1, 20834017.68103584),
    @(1.3, -1, 116738344.7276593, 17872321.7534636, -1, 20834017.68103584),
    @(1.5, -1, 116738344.7276593, 11521769.02845364, -1, 20834017.68103584),
    @(1.8, -1, 116738344.7276593, 11302843.50130033, -1, 20834017.68103584),
    @(2.5, -1, 80137455.6064499, 14221414.84596121, -1, 11380037.01593464),
    @(3.2, -1, 90371559.83988744, 14822320.3418464, -1, 12431833.00443093),
    @(3.5, -1, 53260479.12785688, 17305511.7539629, -1, 8061486.0910345),
    @(3.6, -1, 116117562.7339392, 15539820.14622358, -1, 15123385.66730704),
    @(3.73, -1, 116962796.7446434, 11825985.62663229, -1, 16344229.97077986),
    @(3.86, -1, 112412969.9542895, 9821637.440598581, -1, 18381504.65203911),
    @(4.9, -1, 132649852.1243531, 16789610.85197198, -1, 17332358.19153716),
    @(6.3, -1, 133421549.2734301, 14595190.84249189, -1, 19107628.50224325),
    @(7.6, -1, 121054749.0576133, 16555300.36441313, -1, 17792521.86196001),
    @(8.4, -1, 123463669.8591429, 15567176.69019014, -1, 19631831.12930312),
    @(9.1, -1, 118821528.4979442, 15473100.96034378, -1, 19887520.68009411),
    @(9.300000000000001, -1, 128566464.3800524, 16958663.53070267, -1, 17826434.91396135),
    @(9.6, -1, 119972532.8467586, 16365996.77554328, -1, 19805851.48055682)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# New rows 13-21 need the same formatting (bold, centered, bordered) as the
# rest of column A, since they fall outside the original A1:F12 range.
$ws.Range("A2").Copy()
$ws.Range("A13:A21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()
